{"js": "// Remove the errant trailing semicolon from the \"_h2o_keep_element;\"\n// marker text that appears in paragraphs styled \"Node End\" (export of\n// export/node.html). Paragraphs using other styles (e.g. \"Node Start\")\n// already read \"_h2o_keep_element\" with no semicolon and must stay as-is.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst TARGET_STYLE = \"Node End\";\nconst OLD_TEXT = \"_h2o_keep_element;\";\nconst NEW_TEXT = \"_h2o_keep_element\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.style === TARGET_STYLE && paragraph.text === OLD_TEXT) {\n    paragraph.insertText(NEW_TEXT, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the errant trailing semicolon from the \"_h2o_keep_element;\"\n# marker text that appears in paragraphs styled \"Node End\" (export of\n# export/node.html). Paragraphs using other styles (e.g. \"Node Start\")\n# already read \"_h2o_keep_element\" with no semicolon and must stay as-is.\n$d = $word.ActiveDocument\n\n$targetStyle = \"Node End\"\n$oldText = \"_h2o_keep_element;\"\n$newText = \"_h2o_keep_element\"\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    $coreText = $p.Range.Text.TrimEnd(\"`r\", \"`f\")\n    if ($styleName -eq $targetStyle -and $coreText -eq $oldText) {\n        $p.Range.Text = $newText\n    }\n}\n"}
